# The sheet originally has data in columns A:F where column A (style s1,
# values 9/15) is a leftover column that should be removed. Removing the
# entire column A and shifting everything left reproduces the target
# layout (A:E) exactly: old B->A, C->B, D->C, E->D, F->E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()
